# CIERRE 7 DIC 2021
# - Switch the active/selected tab from "ARQUITECTO" to "VALES DE INSENTIVOS"
# - Update the incentive-month text from OCTUBRE to NOVIEMBRE
# - TODAY() volatile formulas recalc automatically on save

$wb = $excel.ActiveWorkbook

$wsArquitecto = $wb.Worksheets.Item("ARQUITECTO        ")
$wsVales      = $wb.Worksheets.Item("VALES DE INSENTIVOS")

# Update the incentive period text (shared string used in VALES DE INSENTIVOS!A4)
$wsVales.Range("A4").Value = "PAGO DE INCENTIVO DEL MES DE NOVIEMBRE  2021"

# Make "VALES DE INSENTIVOS" the active/selected sheet (was "ARQUITECTO")
$wsVales.Activate()
